# Update Name of Algo
# Applies updated KNN imputation result values to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "E3"  = 16.195
    "B12" = 5.492999999999999
    "E14" = 17.237
    "E26" = 16.117
    "B27" = 5.119000000000001
    "E31" = 16.139
    "B32" = 6.284000000000001
    "E35" = 16.63
    "B36" = 7.987
    "E37" = 16.737
    "B38" = 5.287000000000001
    "E45" = 16.749
    "B46" = 6.165
    "E52" = 16.768
    "B54" = 5.072
    "B55" = 4.684
    "B56" = 4.737
    "E57" = 16.529
    "B67" = 4.989
    "B69" = 5.002
    "B72" = 5.380999999999999
    "E81" = 16.268
    "B83" = 5.883
    "E83" = 16.973
    "B86" = 5.016
    "B91" = 6.026999999999999
    "B93" = 5.380999999999999
    "B99" = 5.217000000000001
    "E100" = 16.568
    "E102" = 16.279
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
